# Fill in last week's hours (rows 42-55, columns C=Start time, D=Lunch break, E=End Time)
# Column F (Hours) already has the shared formula =E-C-D, so it recalculates automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Days with no hours logged yet this week still get explicit zeros.
$ws.Range("C42:E49").Value = 0

# Tuesday 45462 (row 50): 11:00 start, no lunch break, 16:00 end
$ws.Cells.Item(50, 3).Value = 11
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(50, 5).Value = 16

# Wednesday 45463 (row 51): 10:00 start, 0.25 lunch, 17:30 end
$ws.Cells.Item(51, 3).Value = 10
$ws.Cells.Item(51, 4).Value = 0.25
$ws.Cells.Item(51, 5).Value = 17.5

# Thursday 45464 (row 52): 10:00 start, 0.25 lunch, 18:00 end
$ws.Cells.Item(52, 3).Value = 10
$ws.Cells.Item(52, 4).Value = 0.25
$ws.Cells.Item(52, 5).Value = 18

# Sunday 45467 (row 55): only a start time of 8.5 logged so far
$ws.Cells.Item(55, 3).Value = 8.5

# Leave the cursor where the author left off
$ws.Range("D56").Select()
